# "Generate Report for Handoff"
#
# The localization-status report previously showed the .md source file's
# handoff as failed, with no handoff file/date recorded. A handoff run has
# now produced real .xlf target files for zh-cn and de-de, so the report is
# updated to reflect that: status flips to "Ready for handoff", the new
# handoff file name + datetime are recorded, and the .localization-config
# row's reason changes from "Ignored" to "Include".

$wb = $excel.ActiveWorkbook

$commitPath = "https://github.com/OpenLocalizationTest/oltest/blob/f95aa65571649722770d9a6ad296f3a68a164924"
$mdName = "66f80cc9-80c1-471b-b042-8610d0a171d2.md"
$cfgName = ".localization-config"

# ---------------------------------------------------------------------
# Overview sheet: per-locale status column for the .md source row.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# Per-locale detail sheets (zh-cn / de-de): record the new handoff file,
# its datetime, and flip the dependency row's reason to "Include".
# ---------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; XlfName = "66f80cc9-80c1-471b-b042-8610d0a171d2.4d7e36a84b2a77b3b8d88cc9f745aa27d84fad08.zh-cn.xlf"; HandoffDate = "2016-02-17 04:44:50" },
    @{ Sheet = "de-de"; XlfName = "66f80cc9-80c1-471b-b042-8610d0a171d2.4d7e36a84b2a77b3b8d88cc9f745aa27d84fad08.de-de.xlf"; HandoffDate = "2016-02-17 04:44:59" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Status + latest handoff datetime for the .md source row.
    $ws.Range("B2").Value = "Ready for handoff"
    $ws.Range("D2").Value = $loc.HandoffDate

    # Dependency row's handoff reason: no longer ignored, now included.
    $ws.Range("H2").Value = "Include"

    # Rebuild the hyperlinks so the new "Latest Handoff File" link lands
    # between the existing two (md source, then new xlf, then config).
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "$commitPath/e2e/$mdName", "", "", $mdName)
    $ws.Hyperlinks.Add($ws.Range("C2"), "$commitPath/e2e/$($loc.XlfName)", "", "", $loc.XlfName)
    $ws.Hyperlinks.Add($ws.Range("A3"), "$commitPath/$cfgName", "", "", $cfgName)
}

Write-Output "Report regenerated for handoff."
